$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: drop the stray empty cells that only held an empty inline
#     string (C13, E13, F13, G13, H13). D13 keeps its email value.
$ws.Range("C13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()

# --- Row 14: brand new request row appended below the existing data.
$ws.Range("A14").Value = "2025-12-16 02:31:30 UTC"
$ws.Range("B14").Value = "sujay jirapure"
$ws.Range("C14").Value = "KGN Solar"
$ws.Range("D14").Value = "jirapuresujay@gmail.com"
$ws.Range("E14").Value = "IN"

# Dial code / phone columns start with a literal "+" followed only by
# digits, which Excel's input parser would otherwise coerce to a number.
# Format them as text first so the leading "+" (and the full digit
# string) is preserved verbatim.
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "+91"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "+9161408591185"

$ws.Range("H14").Value = "50 Broughton Road"
# I14 has no message text for this request - leave it blank.
$ws.Range("J14").Value = "192.168.1.54"
